# Update the USD Amount figure in cell T2 (SheetName1) to reflect the
# latest daily update (52856 -> 53531), matching the new file's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

$ws.Range("T2").Value = 53531
